$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is organized as small "groups" of rows, each preceded by a
# standalone "group header" row that only has a label in column A (the
# thread-size group, e.g. "4-40", "6-32", ...) and the material in column M.
# This edit removes those standalone header rows and instead fills their
# label down into column L ("thread_size") for every data row that belongs
# to the group, then also turns row 1 (a text header row) into a plain
# numeric column-index row (0..12), moving the original text headers down
# to row 2 (with the old L/M header text cleared).

# Process groups from the bottom of the sheet upward so that row numbers
# for groups still to be processed are not disturbed by row deletions.
$groups = @(
    @{ Header = 49; Start = 50; End = 52 },
    @{ Header = 31; Start = 32; End = 38 },
    @{ Header = 24; Start = 25; End = 30 },
    @{ Header = 15; Start = 16; End = 23 },
    @{ Header = 8;  Start = 9;  End = 14 },
    @{ Header = 2;  Start = 3;  End = 7 }
)

# Row 1 currently holds the real text column headers (Lg., Threading, ...).
# Capture that text now (columns A..K only -- L/M header text is dropped)
# so it can be written back into row 2 once the header rows above have
# been removed.
$headerTexts = @()
for ($c = 1; $c -le 11; $c++) {
    $headerTexts += , ($ws.Cells.Item(1, $c).Value())
}

foreach ($g in $groups) {
    $headerRow = $g.Header
    $label = $ws.Cells.Item($headerRow, 1).Value()

    for ($r = $g.Start; $r -le $g.End; $r++) {
        $ws.Cells.Item($r, 12).Value = $label
    }

    $ws.Range("A" + $headerRow + ":M" + $headerRow).EntireRow.Delete()
}

# Move the original text headers down to row 2 (L2/M2 stay blank), and
# replace row 1 with plain numeric column indices 0..12. A new row must be
# inserted for the headers -- row 2 is currently the first real data row,
# and simply overwriting it would clobber that data instead of pushing it
# down. The inserted row inherits row 1's bold/centered formatting, but
# row 2 should end up unstyled like the rest of the data rows, so strip
# that back off again.
$ws.Rows(2).Insert()
$ws.Range("A2:M2").ClearFormats()
for ($c = 1; $c -le 11; $c++) {
    $ws.Cells.Item(2, $c).Value = $headerTexts[$c - 1]
}
$ws.Cells.Item(2, 12).Value = ""
$ws.Cells.Item(2, 13).Value = ""

for ($c = 1; $c -le 13; $c++) {
    $ws.Cells.Item(1, $c).Value = $c - 1
}
